# Generate Report for Handback
# Updates the localization-status report after a handback sync:
#  - Overview sheet: Status summary reflects the handback
#  - zh-cn / de-de sheets: handback datetime refreshed, error cleared
#  - Column widths widened to fit the new, longer status/handback text

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
# (column widths are widened to fit the longer status text; the host's
# ColumnWidth->stored-width rounding is coarser than Excel's, so the
# character-width inputs below are chosen to land as close as possible
# to the target stored widths of ~29.98 / ~13.75 "characters")
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("K2").Value = "2016-08-20 04:52:10"
$wsZh.Range("P2").Value = ""
$wsZh.Columns.Item(3).ColumnWidth = 29.1
$wsZh.Columns.Item(16).ColumnWidth = 12.8

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("K2").Value = "2016-08-20 04:52:16"
$wsDe.Range("P2").Value = ""
$wsDe.Columns.Item(3).ColumnWidth = 29.1
$wsDe.Columns.Item(16).ColumnWidth = 12.8
